$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3993.625
$ws.Range("I74").Value = 3449
$ws.Range("J74").Value = 4071.4285
$ws.Range("K74").Value = 3449
$ws.Range("L74").Value = 4071.4285
$ws.Range("M74").Value = -2513
$ws.Range("N74").Value = -5943.4285

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3993.625
$ws.Range("I77").Value = 3449
$ws.Range("J77").Value = 4071.4285
$ws.Range("K77").Value = 17245
$ws.Range("L77").Value = 20357.1425
$ws.Range("M77").Value = -12565
$ws.Range("N77").Value = -29717.1425

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4087.7542
$ws.Range("I138").Value = 1544.7567
$ws.Range("J138").Value = 8008.2085
$ws.Range("K138").Value = 4634.2701
$ws.Range("L138").Value = 24024.6255
$ws.Range("M138").Value = 505.7299000000003
$ws.Range("N138").Value = -34304.62549999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 857482.5600000001
$ws.Range("I141").Value = 2023.75
$ws.Range("K141").Value = 6071.25
$ws.Range("M141").Value = -891.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5209603.5
$ws.Range("I2").Value = 8929527
$ws.Range("J2").Value = 1709.95
$ws.Range("K2").Value = 8929527
$ws.Range("L2").Value = 1709.95
$ws.Range("M2").Value = -8929414
$ws.Range("N2").Value = -1935.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3091.05
$ws.Range("I32").Value = 2970.1545
$ws.Range("K32").Value = 2970.1545
$ws.Range("M32").Value = -2683.1545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3108.35
$ws.Range("I61").Value = 2619.7
$ws.Range("J61").Value = 3597
$ws.Range("K61").Value = 2619.7
$ws.Range("L61").Value = 3597
$ws.Range("M61").Value = -2407.7
$ws.Range("N61").Value = -4021

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 5209603.5
$ws.Range("I116").Value = 8929527
$ws.Range("J116").Value = 1709.95
$ws.Range("K116").Value = 8929527
$ws.Range("L116").Value = 1709.95
$ws.Range("M116").Value = -8927233
$ws.Range("N116").Value = -6297.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1770.4166
$ws.Range("I122").Value = 1489.3864
$ws.Range("J122").Value = 2543.25
$ws.Range("K122").Value = 4468.1592
$ws.Range("L122").Value = 7629.75
$ws.Range("M122").Value = -2018.1592
$ws.Range("N122").Value = -12529.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2287.2708
$ws.Range("I132").Value = 1654.8918
$ws.Range("J132").Value = 4414.364
$ws.Range("K132").Value = 4964.6754
$ws.Range("L132").Value = 13243.092
$ws.Range("M132").Value = -2434.6754
$ws.Range("N132").Value = -18303.092

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3108.35
$ws.Range("I136").Value = 2619.7
$ws.Range("J136").Value = 3597
$ws.Range("K136").Value = 7859.099999999999
$ws.Range("L136").Value = 10791
$ws.Range("M136").Value = -5309.099999999999
$ws.Range("N136").Value = -15891

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5209603.5
$ws.Range("I3").Value = 8929527
$ws.Range("J3").Value = 1709.95
$ws.Range("K3").Value = 8929527
$ws.Range("L3").Value = 1709.95
$ws.Range("M3").Value = -8929413
$ws.Range("N3").Value = -1937.95

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1824.25
$ws.Range("I105").Value = 1386.6666
$ws.Range("J105").Value = 3137
$ws.Range("K105").Value = 1386.6666
$ws.Range("L105").Value = 3137
$ws.Range("M105").Value = 360.3334
$ws.Range("N105").Value = -6631

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2590.0728
$ws.Range("I134").Value = 2570.875
$ws.Range("K134").Value = 7712.625
$ws.Range("M134").Value = -5177.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 3366.0588
$ws.Range("I70").Value = 1844.8
$ws.Range("J70").Value = 3999.9167
$ws.Range("K70").Value = 5534.4
$ws.Range("L70").Value = 11999.7501
$ws.Range("M70").Value = -5219.4
$ws.Range("N70").Value = -12629.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 3366.0588
$ws.Range("I73").Value = 1844.8
$ws.Range("J73").Value = 3999.9167
$ws.Range("K73").Value = 5534.4
$ws.Range("L73").Value = 11999.7501
$ws.Range("M73").Value = -4442.4
$ws.Range("N73").Value = -14183.7501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 3000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 3000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1937.2632
$ws.Range("I118").Value = 271.6
$ws.Range("J118").Value = 2532.1428
$ws.Range("K118").Value = 814.8000000000001
$ws.Range("L118").Value = 7596.428400000001
$ws.Range("M118").Value = 428.1999999999999
$ws.Range("N118").Value = -10082.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 24066.738
$ws.Range("I129").Value = 3783.75
$ws.Range("J129").Value = 34884.332
$ws.Range("K129").Value = 11351.25
$ws.Range("L129").Value = 104652.996
$ws.Range("M129").Value = -6351.25
$ws.Range("N129").Value = -114652.996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1272.081
$ws.Range("I131").Value = 2814.4443
$ws.Range("J131").Value = 1058.5231
$ws.Range("K131").Value = 8443.332900000001
$ws.Range("L131").Value = 3175.5693
$ws.Range("M131").Value = -3403.332900000001
$ws.Range("N131").Value = -13255.5693

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3612.4075
$ws.Range("I122").Value = 3096
$ws.Range("J122").Value = 4363.5454
$ws.Range("K122").Value = 9288
$ws.Range("L122").Value = 13090.6362
$ws.Range("M122").Value = -6838
$ws.Range("N122").Value = -17990.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4041.4595
$ws.Range("I132").Value = 4631.5293
$ws.Range("K132").Value = 13894.5879
$ws.Range("M132").Value = -11364.5879

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1293.2162
$ws.Range("I46").Value = 943.9394
$ws.Range("J46").Value = 4174.75
$ws.Range("K46").Value = 943.9394
$ws.Range("L46").Value = 4174.75
$ws.Range("M46").Value = -755.9394
$ws.Range("N46").Value = -4550.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5422.5127
$ws.Range("I132").Value = 1867.9445
$ws.Range("J132").Value = 8469.286
$ws.Range("K132").Value = 5603.833500000001
$ws.Range("L132").Value = 25407.858
$ws.Range("M132").Value = -3073.833500000001
$ws.Range("N132").Value = -30467.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 35528.9
$ws.Range("I126").Value = 41462.68
$ws.Range("J126").Value = 5860
$ws.Range("K126").Value = 124388.04
$ws.Range("L126").Value = 17580
$ws.Range("M126").Value = -121918.04
$ws.Range("N126").Value = -22520
